# Update "Training Dashboard" sheet: decrement "PERIOD TO EXPIRE" (col H) by 1
# and bump "LAST UPDATE" (col I) from 03-Nov-2025 to 04-Nov-2025 for rows 3-26.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Build the new date as a literal text value in an unused scratch cell, well
# outside the worksheet's used range. Pre-formatting it as Text ("@") before
# assigning the date-looking string keeps Excel from auto-converting it into a
# date serial number, the way it would if we wrote the string straight into
# the target cell's .Value.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.NumberFormat = "@"
$scratch.Value = "04-Nov-2025"
$scratch.Copy()

for ($row = 3; $row -le 26; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H - PERIOD TO EXPIRE
    $iCell = $ws.Cells.Item($row, 9)   # column I - LAST UPDATE

    # Decrement the numeric "period to expire" value.
    $hCell.Value = $hCell.Value2 - 1

    # Paste only the scratch cell's value (the literal text "04-Nov-2025")
    # into column I, leaving the target cell's existing formatting/style
    # untouched.
    $iCell.PasteSpecial(-4163)  # xlPasteValues
}

# Clean up: remove the scratch cell and any leftover clipboard marquee.
$scratch.Clear()
$excel.CutCopyMode = $false
